$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for columns G and H, reusing the header style from A1:F1
$ws.Range("G1").Value = "num_samples"
$ws.Range("H1").Value = "fractional_uncertainty"

$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data for num_samples (G) and fractional_uncertainty (H) per row
$numSamples = @{
    2 = 1000
    3 = 1000
    4 = 1000
    5 = 1000
    6 = 1000
    7 = 1000
    8 = 1000
    9 = 1000
    10 = 1000
    11 = 1000
    12 = 993
    13 = 1000
    14 = 988
    15 = 1000
    16 = 999
    17 = 1000
    18 = 995
    19 = 1000
}

$fracUncertainty = @{
    2 = 0.05465967198530568
    3 = 0.1517863286558169
    4 = 0.03076134861937979
    5 = 0.1002143141730702
    6 = 0.03284310601598622
    7 = 0.07607541391771377
    8 = 0.03160184766307855
    9 = 0.05038480137631279
    10 = 0.03108241073896207
    11 = 0.05481633040125691
    12 = 0.1057280547432848
    13 = 0.03084771414731377
    14 = 0.1337567131970655
    15 = 0.05682149308062914
    16 = 0.127338086543459
    17 = 0.06317416027103828
    18 = 0.1301786259478235
    19 = 0.05576595618402488
}

for ($row = 2; $row -le 19; $row++) {
    $ws.Cells.Item($row, 7).Value = $numSamples[$row]
    $ws.Cells.Item($row, 8).Value = $fracUncertainty[$row]
}
